# Apply the recorded edits to the "Test-Cases" sheet of the Camp_2 workbook.
#
# Summary of the change (from the author's re-upload of the workbook):
#   - Row 10's "Approved/Rejected" cell (I10) was changed from "Approved" to "Rejected".
#   - Row 10's "ReasonToReject" cell (J10) was filled in with "checkk".
#   - The active selection on the sheet moved from H17 to H13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Mark the test case in row 10 as Rejected and record the reason.
$ws.Range("I10").Value = "Rejected"
$ws.Range("J10").Value = "checkk"

# Update the sheet's current selection to match the saved view state.
$ws.Range("H13").Select()
